$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix total marks error: update "Marking" row (row 11) and "Total" row (row 12)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "92 / 112"
